# Apply the changes described by the diff:
#  1. Rename the worksheet from "CopperA-HW15.xpc" to "CopperA"
#  2. Nudge a handful of existing cells in rows 13 and 15 by one ULP
#     (values recomputed by the new Gaussian-Quadrature averaging routine)
#  3. Append a new data row (row 16) for HKL index 14 / "HexGrid-60degTilt5degRes"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet
$ws.Name = "CopperA"

# 2. Tiny last-digit precision corrections on row 13
$ws.Range("D13").Value = 0.9934107449639451
$ws.Range("E13").Value = 0.9975992194798802
$ws.Range("H13").Value = 0.9934107449639451
$ws.Range("I13").Value = 0.9939544707871102
$ws.Range("M13").Value = 0.9886723630412961

# ... and on row 15
$ws.Range("C15").Value = 0.9257927643229713
$ws.Range("G15").Value = 0.9257927643229713
$ws.Range("M15").Value = 0.9256325400147358

# 3. Append new row 16, matching the layout/style of the preceding rows.
#    Copy A15's style (bold, bordered, centered) down onto A16 first.
$ws.Range("A15").Copy($ws.Range("A16"))

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 1.160794165004746
$ws.Range("D16").Value = 0.8636495835432553
$ws.Range("E16").Value = 1.19234943773164
$ws.Range("F16").Value = 0.8687421159298508
$ws.Range("G16").Value = 1.160794165004746
$ws.Range("H16").Value = 0.8636495835432553
$ws.Range("I16").Value = 1.075637167713381
$ws.Range("J16").Value = 0.9958781954106263
$ws.Range("K16").Value = 0.9521629437015247
$ws.Range("L16").Value = 0.8499167015456121
$ws.Range("M16").Value = 1.160794165004746
$ws.Range("N16").Value = 1.027999510637448
$ws.Range("O16").Value = 1.021383825552373
$ws.Range("P16").Value = 0.9948912888225796

Write-Host "Applied sheet rename, precision fixes, and new row 16."
